# Handback status report regeneration:
#   - old file pair  01bbae0a-7cbe-451a-9851-9b39a347bbca / 212c757b-9df1-46f1-8d2b-b5f317d9d6a2
#     is replaced by the new pair 8809487c-5638-45af-a084-8eaf3a4785cb / ffff56d9cf63-a69d-471b-84ea-a93c079cdc90
#   - handoff/handback timestamps and hashes move forward to the new run
$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "8809487c-5638-45af-a084-8eaf3a4785cb.md"
$ws1.Range("B2").Value = "e2e\8809487c-5638-45af-a084-8eaf3a4785cb.md"
$ws1.Range("G2").Value = "2016-09-07 11:28:52"

$ws1.Range("A3").Value = "ffff56d9cf63-a69d-471b-84ea-a93c079cdc90.md"
$ws1.Range("B3").Value = "e2e\ffff56d9cf63-a69d-471b-84ea-a93c079cdc90.md"
$ws1.Range("G3").Value = "2016-09-07 11:28:52"

# ---- zh-cn sheet ------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "8809487c-5638-45af-a084-8eaf3a4785cb.md"
$ws2.Range("G2").Value = "8809487c-5638-45af-a084-8eaf3a4785cb.831db3405c1d1424eb912caf903eb320fb28d6be.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-07 11:28:42"
$ws2.Range("I2").Value = "8809487c-5638-45af-a084-8eaf3a4785cb.md"
$ws2.Range("J2").Value = "8809487c-5638-45af-a084-8eaf3a4785cb.831db3405c1d1424eb912caf903eb320fb28d6be.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-07 11:29:20"

$ws2.Range("A3").Value = "ffff56d9cf63-a69d-471b-84ea-a93c079cdc90.md"
$ws2.Range("G3").Value = "8809487c-5638-45af-a084-8eaf3a4785cb.831db3405c1d1424eb912caf903eb320fb28d6be.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-09-07 11:28:42"
$ws2.Range("I3").Value = "ffff56d9cf63-a69d-471b-84ea-a93c079cdc90.md"
$ws2.Range("J3").Value = "8809487c-5638-45af-a084-8eaf3a4785cb.831db3405c1d1424eb912caf903eb320fb28d6be.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-09-07 11:29:20"

# ---- de-de sheet ------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "8809487c-5638-45af-a084-8eaf3a4785cb.md"
$ws3.Range("G2").Value = "8809487c-5638-45af-a084-8eaf3a4785cb.831db3405c1d1424eb912caf903eb320fb28d6be.de-de.xlf"
$ws3.Range("H2").Value = "2016-09-07 11:28:52"
$ws3.Range("I2").Value = "8809487c-5638-45af-a084-8eaf3a4785cb.md"
$ws3.Range("J2").Value = "8809487c-5638-45af-a084-8eaf3a4785cb.831db3405c1d1424eb912caf903eb320fb28d6be.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-07 11:29:29"

$ws3.Range("A3").Value = "ffff56d9cf63-a69d-471b-84ea-a93c079cdc90.md"
$ws3.Range("G3").Value = "8809487c-5638-45af-a084-8eaf3a4785cb.831db3405c1d1424eb912caf903eb320fb28d6be.de-de.xlf"
$ws3.Range("H3").Value = "2016-09-07 11:28:52"
$ws3.Range("I3").Value = "ffff56d9cf63-a69d-471b-84ea-a93c079cdc90.md"
$ws3.Range("J3").Value = "8809487c-5638-45af-a084-8eaf3a4785cb.831db3405c1d1424eb912caf903eb320fb28d6be.de-de.xlf"
$ws3.Range("K3").Value = "2016-09-07 11:29:29"

Write-Host "Handback status report regenerated"
